$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A7").Value = "Data Entry"
$ws.Range("B7").Value = "Auxiliar fuel"
$ws.Range("C7").Value = "number"
$ws.Range("E7").Value = "yes"

$ws.Range("E6").Copy()
$ws.Range("E7").PasteSpecial(-4122)

$ws.Range("E8").Select()
